$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter constant (used by existing cell style, alignment horizontal=center vertical=center)
$xlCenter = -4108

function Set-CellValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    # Match the style (centered horizontal/vertical) used throughout this table
    $rng.HorizontalAlignment = $xlCenter
    $rng.VerticalAlignment = $xlCenter
}

# Row 1 - headers
Set-CellValue "B1" "Yu Qiao"
Set-CellValue "C1" "Luc Van Gool"
Set-CellValue "D1" "Lei Zhang"

# Row 2 - 2022
Set-CellValue "B2" 4
Set-CellValue "C2" 24
Set-CellValue "D2" 15

# Row 3 - 2023
Set-CellValue "B3" 25
Set-CellValue "C3" 15
Set-CellValue "D3" 21

# Row 4 - 2024
Set-CellValue "B4" 26
Set-CellValue "C4" 15
Set-CellValue "D4" 12

# Row 5 - Total
Set-CellValue "B5" 55
Set-CellValue "C5" 54
Set-CellValue "D5" 48
